$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 0.05435008928179741
$ws.Cells.Item(2,2).Value = 0.9844599366188049
$ws.Cells.Item(2,3).Value = 0.01748823188245296
$ws.Cells.Item(2,4).Value = 0.9967902898788452
$ws.Cells.Item(3,1).Value = 0.008115260861814022
$ws.Cells.Item(3,2).Value = 0.9985508918762207
$ws.Cells.Item(3,3).Value = 0.01094591990113258
$ws.Cells.Item(3,4).Value = 0.9973252415657043
$ws.Cells.Item(4,1).Value = 0.004289192613214254
$ws.Cells.Item(4,2).Value = 0.9987987279891968
$ws.Cells.Item(4,3).Value = 0.001903848722577095
$ws.Cells.Item(4,4).Value = 0.9992867112159729
$ws.Cells.Item(5,1).Value = 0.00175549800042063
$ws.Cells.Item(5,2).Value = 0.9995042681694031
$ws.Cells.Item(5,3).Value = 0.001479708007536829
$ws.Cells.Item(5,4).Value = 0.9996433854103088
$ws.Cells.Item(6,1).Value = 0.001832536770962179
$ws.Cells.Item(6,2).Value = 0.9995423555374146
$ws.Cells.Item(6,3).Value = 0.001561190700158477
$ws.Cells.Item(6,4).Value = 0.9996433854103088
$ws.Cells.Item(7,1).Value = 0.001302260672673583
$ws.Cells.Item(7,2).Value = 0.999733030796051
$ws.Cells.Item(7,3).Value = 0.000780583533924073
$ws.Cells.Item(7,4).Value = 0.999821662902832
$ws.Cells.Item(8,1).Value = 0.0008436710340902209
$ws.Cells.Item(8,2).Value = 0.999733030796051
$ws.Cells.Item(8,3).Value = 0.001642719144001603
$ws.Cells.Item(8,4).Value = 0.9994650483131409
$ws.Cells.Item(9,1).Value = 0.001200058031827211
$ws.Cells.Item(9,2).Value = 0.9996758699417114
$ws.Cells.Item(9,3).Value = 0.000322716950904578
$ws.Cells.Item(9,4).Value = 0.999821662902832
$ws.Cells.Item(10,1).Value = 0.0006467446219176054
$ws.Cells.Item(10,2).Value = 0.9997902512550354
$ws.Cells.Item(10,3).Value = 0.0002208898804383352
$ws.Cells.Item(10,4).Value = 0.999821662902832
$ws.Cells.Item(11,1).Value = 0.0004124153929296881
$ws.Cells.Item(11,2).Value = 0.9998474717140198
$ws.Cells.Item(11,3).Value = 0.00001239955417986494
$ws.Cells.Item(11,4).Value = 1
$ws.Cells.Item(12,1).Value = 0.0004468293336685747
$ws.Cells.Item(12,2).Value = 0.999885618686676
$ws.Cells.Item(12,3).Value = 0.0001524301915196702
$ws.Cells.Item(12,4).Value = 0.999821662902832
$ws.Cells.Item(13,1).Value = 0.0004008092801086605
$ws.Cells.Item(13,2).Value = 0.9998093247413635
$ws.Cells.Item(13,3).Value = 0.000004625871952157468
$ws.Cells.Item(13,4).Value = 1
$ws.Cells.Item(14,1).Value = 0.0004431596025824547
$ws.Cells.Item(14,2).Value = 0.999885618686676
$ws.Cells.Item(14,3).Value = 0.0000391287831007503
$ws.Cells.Item(14,4).Value = 1
$ws.Cells.Item(15,1).Value = 0.001117561594583094
$ws.Cells.Item(15,2).Value = 0.9997902512550354
$ws.Cells.Item(15,3).Value = 0.000002446358848828822
$ws.Cells.Item(15,4).Value = 1
$ws.Cells.Item(16,1).Value = 0.0002002900291699916
$ws.Cells.Item(16,2).Value = 0.999885618686676
$ws.Cells.Item(16,3).Value = 0.00001880382478702813
$ws.Cells.Item(16,4).Value = 1
$ws.Cells.Item(17,1).Value = 0.0006398882833309472
$ws.Cells.Item(17,2).Value = 0.9997902512550354
$ws.Cells.Item(17,3).Value = 0.0003039447474293411
$ws.Cells.Item(17,4).Value = 0.9996433854103088
$ws.Cells.Item(18,1).Value = 0.0006364987348206341
$ws.Cells.Item(18,2).Value = 0.9998474717140198
$ws.Cells.Item(18,3).Value = 0.00000466250594399753
$ws.Cells.Item(18,4).Value = 1
$ws.Cells.Item(19,1).Value = 0.0001642543356865644
$ws.Cells.Item(19,2).Value = 0.9999809265136719
$ws.Cells.Item(19,3).Value = 0.00001440150026610354
$ws.Cells.Item(19,4).Value = 1
$ws.Cells.Item(20,1).Value = 0.0004743355966638774
$ws.Cells.Item(20,2).Value = 0.9998665452003479
$ws.Cells.Item(20,3).Value = 0.000001666354478402354
$ws.Cells.Item(20,4).Value = 1
$ws.Cells.Item(21,1).Value = 0.0003056821005884558
$ws.Cells.Item(21,2).Value = 0.9999427795410156
$ws.Cells.Item(21,3).Value = 0.0000140697784445365
$ws.Cells.Item(21,4).Value = 1
$ws.Cells.Item(22,1).Value = 0.0002064280706690624
$ws.Cells.Item(22,2).Value = 0.9999046325683594
$ws.Cells.Item(22,3).Value = 0.000002480159082551836
$ws.Cells.Item(22,4).Value = 1
$ws.Cells.Item(23,1).Value = 0.000159033399540931
$ws.Cells.Item(23,2).Value = 0.9999427795410156
$ws.Cells.Item(23,3).Value = 0.00000335293952957727
$ws.Cells.Item(23,4).Value = 1
$ws.Cells.Item(24,1).Value = 0.0005081337876617908
$ws.Cells.Item(24,2).Value = 0.9998665452003479
$ws.Cells.Item(24,3).Value = 0.0002165497135138139
$ws.Cells.Item(24,4).Value = 0.999821662902832
$ws.Cells.Item(25,1).Value = 0.0002997218689415604
$ws.Cells.Item(25,2).Value = 0.9999237060546875
$ws.Cells.Item(25,3).Value = 0.00001431106647942215
$ws.Cells.Item(25,4).Value = 1
$ws.Cells.Item(26,1).Value = 0.0002636550343595445
$ws.Cells.Item(26,2).Value = 0.9999427795410156
$ws.Cells.Item(26,3).Value = 0.000003478111921140226
$ws.Cells.Item(26,4).Value = 1
$ws.Cells.Item(27,1).Value = 0.00004461240678210743
$ws.Cells.Item(27,2).Value = 0.9999809265136719
$ws.Cells.Item(27,3).Value = 0.00006517051224363968
$ws.Cells.Item(27,4).Value = 1
$ws.Cells.Item(28,1).Value = 0.0007274065283127129
$ws.Cells.Item(28,2).Value = 0.9999046325683594
$ws.Cells.Item(28,3).Value = 0.00001134677495429059
$ws.Cells.Item(28,4).Value = 1
$ws.Cells.Item(29,1).Value = 0.0002211982209701091
$ws.Cells.Item(29,2).Value = 0.9999618530273438
$ws.Cells.Item(29,3).Value = 0.001560364267788827
$ws.Cells.Item(29,4).Value = 0.999821662902832
$ws.Cells.Item(30,1).Value = 0.0004530835140030831
$ws.Cells.Item(30,2).Value = 0.9999237060546875
$ws.Cells.Item(30,3).Value = 0.00001826351399358828
$ws.Cells.Item(30,4).Value = 1
$ws.Cells.Item(31,1).Value = 0.0001113592588808388
$ws.Cells.Item(31,2).Value = 0.9999618530273438
$ws.Cells.Item(31,3).Value = 0.0003492028918117285
$ws.Cells.Item(31,4).Value = 0.999821662902832
$ws.Cells.Item(32,1).Value = 0.0003438794519752264
$ws.Cells.Item(32,2).Value = 0.9999427795410156
$ws.Cells.Item(32,3).Value = 0.00002056075572909322
$ws.Cells.Item(32,4).Value = 1
$ws.Cells.Item(33,1).Value = 0.0002745004603639245
$ws.Cells.Item(33,2).Value = 0.9998665452003479
$ws.Cells.Item(33,3).Value = 0.000001885717097138695
$ws.Cells.Item(33,4).Value = 1
$ws.Cells.Item(34,1).Value = 0.0005310648120939732
$ws.Cells.Item(34,2).Value = 0.9999237060546875
$ws.Cells.Item(34,3).Value = 0.000006786219273635652
$ws.Cells.Item(34,4).Value = 1
$ws.Cells.Item(35,1).Value = 0.0004477001202758402
$ws.Cells.Item(35,2).Value = 0.999885618686676
$ws.Cells.Item(35,3).Value = 0.000002977098574774573
$ws.Cells.Item(35,4).Value = 1
$ws.Cells.Item(36,1).Value = 0.00003265546911279671
$ws.Cells.Item(36,2).Value = 1
$ws.Cells.Item(36,3).Value = 0.00001813925882743206
$ws.Cells.Item(36,4).Value = 1
$ws.Cells.Item(37,1).Value = 0.00002224328272859566
$ws.Cells.Item(37,2).Value = 1
$ws.Cells.Item(37,3).Value = 0.0000003085166611072054
$ws.Cells.Item(37,4).Value = 1
$ws.Cells.Item(38,1).Value = 0.0002350623981328681
$ws.Cells.Item(38,2).Value = 0.9999046325683594
$ws.Cells.Item(38,3).Value = 0.000006502752512460575
$ws.Cells.Item(38,4).Value = 1
$ws.Cells.Item(39,1).Value = 0.0004320971493143588
$ws.Cells.Item(39,2).Value = 0.9999237060546875
$ws.Cells.Item(39,3).Value = 0.000007871895832067821
$ws.Cells.Item(39,4).Value = 1
$ws.Cells.Item(40,1).Value = 0.0001643529103603214
$ws.Cells.Item(40,2).Value = 0.9999237060546875
$ws.Cells.Item(40,3).Value = 0.000008416612217843067
$ws.Cells.Item(40,4).Value = 1
$ws.Cells.Item(41,1).Value = 0.0001762539177434519
$ws.Cells.Item(41,2).Value = 0.9999427795410156
$ws.Cells.Item(41,3).Value = 0.00000003116148050708034
$ws.Cells.Item(41,4).Value = 1
$ws.Cells.Item(42,1).Value = 0.00003312867556815036
$ws.Cells.Item(42,2).Value = 0.9999809265136719
$ws.Cells.Item(42,3).Value = 0.00000006412402342448331
$ws.Cells.Item(42,4).Value = 1
$ws.Cells.Item(43,1).Value = 0.0005391839076764882
$ws.Cells.Item(43,2).Value = 0.9997902512550354
$ws.Cells.Item(43,3).Value = 0.00001292183515033685
$ws.Cells.Item(43,4).Value = 1
$ws.Cells.Item(44,1).Value = 0.0002767968690022826
$ws.Cells.Item(44,2).Value = 0.9999618530273438
$ws.Cells.Item(44,3).Value = 0.000003666962356874137
$ws.Cells.Item(44,4).Value = 1
$ws.Cells.Item(45,1).Value = 0.0001841589546529576
$ws.Cells.Item(45,2).Value = 0.9999427795410156
$ws.Cells.Item(45,3).Value = 0.001688022050075233
$ws.Cells.Item(45,4).Value = 0.999821662902832
$ws.Cells.Item(46,1).Value = 0.0002466020523570478
$ws.Cells.Item(46,2).Value = 0.9999046325683594
$ws.Cells.Item(46,3).Value = 0.0000009921159289660864
$ws.Cells.Item(46,4).Value = 1
$ws.Cells.Item(47,1).Value = 0.0000635948235867545
$ws.Cells.Item(47,2).Value = 0.9999809265136719
$ws.Cells.Item(47,3).Value = 0.000001256919176739757
$ws.Cells.Item(47,4).Value = 1
$ws.Cells.Item(48,1).Value = 0.00001011686163110426
$ws.Cells.Item(48,2).Value = 1
$ws.Cells.Item(48,3).Value = 0.000003656909484561766
$ws.Cells.Item(48,4).Value = 1
$ws.Cells.Item(49,1).Value = 0.000664616993162781
$ws.Cells.Item(49,2).Value = 0.999885618686676
$ws.Cells.Item(49,3).Value = 0.00000007917800104451089
$ws.Cells.Item(49,4).Value = 1
$ws.Cells.Item(50,1).Value = 0.0007033126894384623
$ws.Cells.Item(50,2).Value = 0.9998665452003479
$ws.Cells.Item(50,3).Value = 0.00002501680683053564
$ws.Cells.Item(50,4).Value = 1
$ws.Cells.Item(51,1).Value = 0.00001468202026444487
$ws.Cells.Item(51,2).Value = 1
$ws.Cells.Item(51,3).Value = 0.00001690940916887484
$ws.Cells.Item(51,4).Value = 1
